# Saldo_guide.xlsx update: refresh the daily snapshot
#  - sheet renamed from the 2024-08-21 run to the 2024-08-22 run
#  - every balance date (column G) rolled from 45525 to 45526 (one day later)
#  - a handful of "total" rows (previously negative D + positive E) now show
#    the prior-day negative amount zeroed out, with E/H carrying the new net

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Sheet (tab) name
$ws.Name = "IClientBalance-20240822-093758-"

# 2) Column G: every data row's date bumps by one day (45525 -> 45526)
$ws.Range("G2:G274").Value = 45526

# 3) Rows whose D/E/H values change beyond the date roll.
#    For each: D -> 0, E -> new net, H -> new net (H = D + E)
$rows = @(
    @{ Row = 17;  E = 12941.59 },
    @{ Row = 43;  E = -85.99 },
    @{ Row = 99;  E = -28.05 },
    @{ Row = 101; E = 20069.98 },
    @{ Row = 103; E = 55063.1 },
    @{ Row = 104; E = -70.99 },
    @{ Row = 108; E = -105.08 },
    @{ Row = 132; E = -70.42 },
    @{ Row = 143; E = -265.09 },
    @{ Row = 158; E = 83.56 },
    @{ Row = 173; E = -67.6 },
    @{ Row = 235; E = -3.49 },
    @{ Row = 264; E = -82.9 },
    @{ Row = 265; E = -45.01 },
    @{ Row = 273; E = -31.54 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $e = $item.E
    $ws.Range("D$r").Value = 0
    $ws.Range("E$r").Value = $e
    $ws.Range("H$r").Value = $e
}
